# Weekly update: a new daily price record for Cilantro (Terminal La Palmera
# de La Serena) is inserted at row 55, pushing the existing rows 55-158 down
# to 56-159 (all their data stays the same, just shifted down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 55; Excel shifts rows 55..158 down to 56..159
# and expands the used range / dimension accordingly.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record's data. The columns
# that are constant for every row in this sheet (A,B,C,E,F,G,H,I,N,O,Q,R) are
# copied verbatim from the surrounding rows; D,J,K,L,M,P hold the new record.
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44775
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112040
$ws.Range("G55").Value = "Cilantro"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = 2250
$ws.Range("N55").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O55").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P55").Value = 1500
$ws.Range("Q55").Value = 1.5
$ws.Range("R55").Value = "Hortaliza"

Write-Output "Inserted new row 55 and shifted subsequent rows down."
